# Update the "想去人数" (want-to-go count) values in column F for a handful
# of rows that are duplicated between the "展览" sheet and the "全部类型" sheet.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsAll     = $wb.Worksheets.Item("全部类型")

# 展览 sheet (sheet1) updates
$wsExhibit.Range("F4").Value  = 620
$wsExhibit.Range("F11").Value = 103
$wsExhibit.Range("F15").Value = 34
$wsExhibit.Range("F22").Value = 204
$wsExhibit.Range("F23").Value = 48
$wsExhibit.Range("F26").Value = 199
$wsExhibit.Range("F27").Value = 60

# 全部类型 sheet (sheet4) updates - same events, duplicated rows
$wsAll.Range("F5").Value  = 620
$wsAll.Range("F13").Value = 103
$wsAll.Range("F18").Value = 34
$wsAll.Range("F29").Value = 204
$wsAll.Range("F30").Value = 48
$wsAll.Range("F35").Value = 199
$wsAll.Range("F36").Value = 60
